$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Recorded By" email lists, attendance counts, and computed
# percentage metrics per the 2025-10-08 13:50:17 attendance processing run.

$ws.Range("G2").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G4").Value = "mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G5").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G8").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("L10").Formula = '="53.8%"'
$ws.Range("L10").Copy()
$ws.Range("L10").PasteSpecial(-4163)
$ws.Range("G14").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G15").Value = "dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G17").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("S17").Formula = '="48.8%"'
$ws.Range("S17").Copy()
$ws.Range("S17").PasteSpecial(-4163)
$ws.Range("G18").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("S18").Formula = '="63.4%"'
$ws.Range("S18").Copy()
$ws.Range("S18").PasteSpecial(-4163)
$ws.Range("G19").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G21").Value = "mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G22").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G25").Value = "Alshimaa_khaled@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg"
$ws.Range("G31").Value = "nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G32").Value = "dr_ahmed_mostafa@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G34").Value = "nardine.alfonse@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg"
$ws.Range("G35").Value = "neveen.nashaat@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg"
$ws.Range("G37").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G45").Value = "System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G48").Value = "Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G49").Value = "Aya.alaa-eldein@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("H49").Value = "176/220"
$ws.Range("G51").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G52").Value = "naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G54").Value = "System, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("G62").Value = "System, Rania.a.youssef@med.asu.edu.eg, backup@backdoor.com, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G65").Value = "Aya.alaa-eldein@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G66").Value = "Aya.alaa-eldein@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("H66").Value = "179/225"
$ws.Range("G68").Value = "abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, naema.gomaa@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, eman.samir@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg"
$ws.Range("G69").Value = "naema.gomaa@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, abdullah.elagrody@med.asu.edu.eg, yasmin.m.senosy@med.asu.edu.eg"
$ws.Range("G70").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G72").Value = "mariam.noureldin@med.asu.edu.eg, Omnia.Mohammed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg, Safa.hany@med.asu.edu.eg"
$ws.Range("G73").Value = "mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G77").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G83").Value = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G85").Value = "yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G86").Value = "Monica.Eshak@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, youstina.magdy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marina_atef@med.asu.edu.eg, Salma.hassan@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G87").Value = "Veronia.rafat@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg"
$ws.Range("G89").Value = "mariam.noureldin@med.asu.edu.eg, Dr.hend.fariid@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, aya.saeed@med.asu.edu.eg"
$ws.Range("G92").Value = "mariam.youssif.std@med.asu.edu.eg, Rania.a.youssef@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G98").Value = "user@user.com, nourhanmohamed@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhanhosni@med.asu.edu.eg, Walaa.h.ghanima@med.asu.edu.eg"
$ws.Range("G102").Value = "yasmin.m.senosy@med.asu.edu.eg, Monica.Eshak@med.asu.edu.eg, Aya_hamed@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"
$ws.Range("G106").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G110").Value = "Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G116").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G117").Value = "Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg"
$ws.Range("G119").Value = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G123").Value = "Sara_nabil@med.asu.edu.eg, aml.awwad@med.asu.edu.eg, norhan.mohamed@med.asu.edu.eg, wessam.atef@med.asu.edu.eg"
$ws.Range("G127").Value = "Alshimaa_khaled@med.asu.edu.eg, afnan.fares@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G133").Value = "enas.omran@med.asu.edu.eg, afaf.abdallah@med.asu.edu.eg, nourhan.mostafa@med.asu.edu.eg"
$ws.Range("G134").Value = "Amr-Saeed@med.asu.edu.eg`n, enas.omran@med.asu.edu.eg, user@user.com, takyas@med.asu.edu.eg"
$ws.Range("G136").Value = "yasmin.m.senosy@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg, nardine.alfonse@med.asu.edu.eg, neveen.nashaat@med.asu.edu.eg, marina_atef@med.asu.edu.eg, ola.m.abdelfattah@med.asu.edu.eg, marinasorial@med.asu.edu.eg, Remon.Matta@med.asu.edu.eg, aya.hanafy@med.asu.edu.eg"
$ws.Range("G144").Value = "nourhanmamdouh@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg, Mayarembaby@med.asu.edu.eg, mohamed.saleem@med.asu.edu.eg"
$ws.Range("G150").Value = "afaf.abdallah@med.asu.edu.eg, Youstina.ibrahim@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
$ws.Range("G151").Value = "System, randa.rabea@med.asu.edu.eg, nourhanmohamed@med.asu.edu.eg, marian.samir@med.asu.edu.eg"
